# Apply the "updated 4.0 files and mdl" edit:
#  - About!C1 date bumped (45320 -> 45392, i.e. 2024-01-29 -> 2024-04-10)
#  - MCF sheet: most non-zero capacity-factor inputs raised to 1 (100%)
#  - MCF sheet: active selection moved from E8 to B17

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsMCF   = $wb.Worksheets.Item("MCF")

# --- About sheet: bump the "last updated" date in C1 ---
$wsAbout.Range("C1").Value = 45392

# --- MCF sheet: raise capacity factors to 1 (100%) ---
$wsMCF.Range("B2").Value  = 1   # hard coal
$wsMCF.Range("B3").Value  = 1   # natural gas steam turbine
$wsMCF.Range("B4").Value  = 1   # natural gas combined cycle
$wsMCF.Range("B6").Value  = 1   # hydro
$wsMCF.Range("B10").Value = 1   # biomass
$wsMCF.Range("B11").Value = 1   # geothermal
$wsMCF.Range("B12").Value = 1   # petroleum
$wsMCF.Range("B13").Value = 1   # natural gas peaker
$wsMCF.Range("B14").Value = 1   # lignite
$wsMCF.Range("B16").Value = 1   # crude oil
$wsMCF.Range("B17").Value = 1   # heavy or residual fuel oil
$wsMCF.Range("B18").Value = 1   # municipal solid waste

# The dependent formula cells B19:B22, B24:B25 (=B2, =B4, =B10, =B14) will
# recalc automatically to 1.

# --- MCF sheet: update the active selection to match the saved view ---
$wsMCF.Activate()
$wsMCF.Range("B17").Select()
